# Update the "dSF" column (F) values for several rows, per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F7").Value  = 1
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = -1
$ws.Range("F14").Value = 4
$ws.Range("F22").Value = 1
$ws.Range("F23").Value = -5
$ws.Range("F26").Value = 2
$ws.Range("F27").Value = -3
